$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 93, shifting existing rows 93-203 down to 94-204.
$ws.Rows.Item(93).EntireRow.Insert()

# Populate the newly-inserted row 93 with the new observation.
$ws.Range("A93").Value = 3
$ws.Range("B93").Value = "Femacal de La Calera"
$ws.Range("C93").Value = "Coquimbo"
$ws.Range("D93").Value = 44483
$ws.Range("E93").Value = 5
$ws.Range("F93").Value = 100114013
$ws.Range("G93").Value = "Zanahoria"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 510
$ws.Range("K93").Value = 8000
$ws.Range("L93").Value = 8500
$ws.Range("M93").Value = 8255
$ws.Range("N93").Value = "$/saco 20 kilos"
$ws.Range("O93").Value = "Chillán"
$ws.Range("P93").Value = 413
$ws.Range("Q93").Value = 20
$ws.Range("R93").Value = "Hortaliza"
